$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweaks -------------------------------------------------
# "Completed" -> "Completed / Log"
$ws.Range("B1").Value = "Completed / Log"

# --- Fix typo in existing Learned/Practiced note (row 3) ---------------
$ws.Range("E3").Value = "Jquery. GIT pushing. Manipulating the DOM with JS. CSS stylzing and element hierarchy. HTML Elements."

# --- Center-align the "Hours Spent" value in row 3 ----------------------
$ws.Range("D3").HorizontalAlignment = -4108  # xlCenter

# --- New log entry: row 4 ----------------------------------------------
$ws.Range("A4").Value = 42883
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)          # xlPasteFormats (reuse date style)

$ws.Range("B4").Value = "Added ""bold,"" ""italic,"" and ""underline"" buttons. Do not function properly yet. Researched File Uploading in Jquery."
$ws.Range("C4").Value = "Chris Peterkin"

$ws.Range("D4").Value = 2
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)          # xlPasteFormats (reuse centered style)

$ws.Range("E4").Value = "Manipulating text in textarea to be bold, italic, or underlined. Options for Jquery file upload. AJAX seems like the best option, it updates files on to the server, but without having to render the page again. Called a dynamic process - manipulates the DOM. There are also plugins for File Uploading that are available.   Created a fork in GITHub."

# --- New section header far below: row 19, column B ---------------------
$ws.Range("A1").Copy()
$ws.Range("B19").PasteSpecial(-4122)         # xlPasteFormats (reuse bold+centered header style)
$ws.Range("B19").Value = "To Do"

$ws.Application.CutCopyMode = 0

# --- Widen column E slightly --------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 120.7

# --- Selection / view state ---------------------------------------------
$ws.Activate()
$ws.Range("E4").Select()
